$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts the existing ExpSch1..ExpSch4 rows down
# from 2-5 to 3-6) to hold the new "Atr" audit-trail schedule entry.
$ws.Rows.Item(2).Insert()

# Row Insert() in this engine copies the formatting of the row above (the
# header row), so strip that back off before writing the new row's content.
$ws.Range("A2:J2").ClearFormats()

$ws.Cells.Item(2, 1).Value = "Atr"
$ws.Cells.Item(2, 2).Value = "test"
$ws.Cells.Item(2, 3).Value = "00:03:00"
$ws.Cells.Item(2, 4).Value = "Daily"
$ws.Cells.Item(2, 5).Value = "AuditTrailReport"
$ws.Cells.Item(2, 8).Value = "EC2AMAZ-N8SAHHO\Administrator"
$ws.Cells.Item(2, 10).Value = "02/06/2021 17:24:49"

# "Date To Export" (G) and "Last Changed On" (J) columns carry the
# dd/MM/yyyy HH:mm:ss date/time display format used by the rest of the sheet.
$ws.Cells.Item(2, 7).NumberFormat = "dd/MM/yyyy HH:mm:ss"
$ws.Cells.Item(2, 10).NumberFormat = "dd/MM/yyyy HH:mm:ss"

# Update the "Last Changed On" timestamps for the schedules that were pushed
# down to rows 3-6 (ExpSch1, ExpSch2, ExpSch3, ExpSch4 respectively).
$ws.Cells.Item(3, 10).Value = "03/06/2021 20:53:07"
$ws.Cells.Item(4, 10).Value = "03/06/2021 20:46:18"
$ws.Cells.Item(5, 10).Value = "03/06/2021 20:47:16"
$ws.Cells.Item(6, 10).Value = "03/06/2021 20:48:11"
